$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp banner in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Junio de 2020 a las 10:35"

# Row 6
$ws.Range("B6").Value = 641156
$ws.Range("C6").Value = 6719
$ws.Range("D6").Value = 403430
$ws.Range("E6").Value = 228560
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 93
$ws.Range("H6").Value = 9166

# Row 37
$ws.Range("A37").Value = "Singapur"
$ws.Range("B37").Value = 43661
$ws.Range("C37").Value = 202
$ws.Range("D37").Value = 37508
$ws.Range("E37").Value = 6127
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 26

# Row 38
$ws.Range("A38").Value = "Ucrania"
$ws.Range("B38").Value = 43628
$ws.Range("C38").Value = 646
$ws.Range("D38").Value = 19027
$ws.Range("E38").Value = 23454
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 18
$ws.Range("H38").Value = 1147

# Row 42
$ws.Range("D42").Value = 20897
$ws.Range("E42").Value = 11572

# Row 79
$ws.Range("A79").Value = "El Salvador"
$ws.Range("B79").Value = 6173
$ws.Range("C79").Value = 239
$ws.Range("D79").Value = 3648
$ws.Range("E79").Value = 2361
$ws.Range("F79").Value = 0
$ws.Range("G79").Value = 12
$ws.Range("H79").Value = 164

# Row 80
$ws.Range("A80").Value = "Republica de Macedonia"
$ws.Range("B80").Value = 6080
$ws.Range("C80").Value = 0
$ws.Range("D80").Value = 2315
$ws.Range("E80").Value = 3479
$ws.Range("F80").Value = 0
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 286

# Row 81
$ws.Range("A81").Value = "Kenia"
$ws.Range("B81").Value = 6070
$ws.Range("C81").Value = 0
$ws.Range("D81").Value = 1971
$ws.Range("E81").Value = 3956
$ws.Range("F81").Value = 0
$ws.Range("G81").Value = 0
$ws.Range("H81").Value = 143

# Row 111
$ws.Range("D111").Value = 1678
$ws.Range("E111").Value = 348

# Row 116
$ws.Range("B116").Value = 1816
$ws.Range("C116").Value = 1
$ws.Range("D116").Value = 1512
$ws.Range("E116").Value = 226

# Row 118
$ws.Range("B118").Value = 1665
$ws.Range("C118").Value = 1
$ws.Range("D118").Value = 1464
$ws.Range("E118").Value = 173

# Row 120
$ws.Range("B120").Value = 1585
$ws.Range("C120").Value = 4
$ws.Range("E120").Value = 90

# Row 131
$ws.Range("B131").Value = 1117
$ws.Range("C131").Value = 1
$ws.Range("E131").Value = 155

# Row 193
$ws.Range("A193").Value = "Islas Turcas y Caicos"
$ws.Range("B193").Value = 41
$ws.Range("C193").Value = 13
$ws.Range("D193").Value = 11
$ws.Range("E193").Value = 29
$ws.Range("F193").Value = 0
$ws.Range("G193").Value = 0
$ws.Range("H193").Value = 1

# Row 194
$ws.Range("A194").Value = "Puerto Rico"
$ws.Range("B194").Value = 39
$ws.Range("C194").Value = 0
$ws.Range("D194").Value = 1
$ws.Range("E194").Value = 36
$ws.Range("F194").Value = 0
$ws.Range("G194").Value = 0
$ws.Range("H194").Value = 2

# Row 195
$ws.Range("A195").Value = "Guam"
$ws.Range("B195").Value = 32
$ws.Range("C195").Value = 0
$ws.Range("D195").Value = 0
$ws.Range("E195").Value = 31
$ws.Range("F195").Value = 0
$ws.Range("G195").Value = 0
$ws.Range("H195").Value = 1

# Row 196
$ws.Range("A196").Value = "San Vicente y las Granadinas"
$ws.Range("B196").Value = 29
$ws.Range("C196").Value = 0
$ws.Range("D196").Value = 29
$ws.Range("E196").Value = 0
$ws.Range("F196").Value = 0
$ws.Range("G196").Value = 0
$ws.Range("H196").Value = 0

# Row 205
$ws.Range("A205").Value = "Dominica"

# Row 206
$ws.Range("A206").Value = "Fiyi"
